$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data one column to the right (A:D -> B:E) and one row
# down (1:16 -> 2:17), preserving the original per-column widths on B:E.
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# New header row across B1:E1 (written first so the shared-string table
# gets these entries before the row labels, matching document order).
$ws.Range("B1").Value2 = "Valid"
$ws.Range("C1").Value2 = "T"
$ws.Range("D1").Value2 = "Z"
$ws.Range("E1").Value2 = "p-value"

# Row labels for the new column A (rows 2..17), one per statistic row.
$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "CyclomaticComplexity(CC) & EffortToImplement",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbUniqueOperands & NbUniqueOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "NbOperators & NbOperators",
    "NbOperators & EffortToImplement",
    "VocabularySize & VocabularySize",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & CyclomaticComplexity(CC)",
    "EffortToImplement & NbOperands",
    "EffortToImplement & NbOperators",
    "EffortToImplement & EffortToImplement",
    "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value2 = $labels[$i]
}

# Widen the new column A to fit the labels (closest achievable value given
# the host's column-width quantization).
$ws.Columns.Item(1).ColumnWidth = 53.6
